$d = $word.ActiveDocument

# Remove the last three paragraphs: the two "In attempting..." / "Our first
# endeavors..." paragraphs and the trailing empty paragraph, leaving the
# paragraph ending in "...I have thrived on my thoughts." as the last piece
# of body text before the section properties.

$count = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($count - 2)
$deleteStart = $startPara.Range.Start
$deleteEnd = $d.Content.End

$r = $d.Range($deleteStart, $deleteEnd)
$r.Delete()
